$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day label (shown as "<day>号"), user count, and font used for the
# trailing "号" character's rich-text run (mirrors how Excel recorded
# each keystroke while the data was typed in).
$data = @(
    @(4, 40, "等线"),
    @(5, 25, "宋体"),
    @(6, 47, "等线"),
    @(7, 19, "等线"),
    @(8, 100, "等线"),
    @(9, 0, "等线"),
    @(10, 2, "等线"),
    @(11, 80, "等线"),
    @(12, 18, "等线")
)

$row = 5
foreach ($entry in $data) {
    $day = $entry[0]
    $count = $entry[1]
    $font = $entry[2]
    $label = "$day" + "号"

    $srcA = $ws.Range("A4")
    $srcA.Copy()
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.PasteSpecial(-4122)  # xlPasteFormats
    $cellA.Value = $label

    $digits = "$day".Length
    $chars = $cellA.Characters($digits + 1, 1)
    $chars.Font.Name = $font
    $chars.Font.Size = 10
    $chars.Font.ColorIndex = -4105  # xlColorIndexAutomatic

    $srcB = $ws.Range("B4")
    $srcB.Copy()
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.PasteSpecial(-4122)  # xlPasteFormats
    $cellB.Value = $count

    $row++
}

$ws.Range("B14").Select()
